$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new column for "M.Sargent" right before the "R.Woods" column (I).
    $ws.Range("I1").EntireColumn.Insert()
    $ws.Range("I1").Value = "M.Sargent"
    $ws.Range("I2").Value = "n"

    # Insert a new column for "K.Blanton" right before the "T.Higbee" column,
    # which (after the previous insert shifted things right) now sits at P.
    $ws.Range("P1").EntireColumn.Insert()
    $ws.Range("P1").Value = "K.Blanton"
    $ws.Range("P2").Value = "n"
}
